# Add new author "C. Tonnel\'e" (DIPC) as a new row in the AuthorList sheet.
# The new entry is inserted at row 104 (pushing the existing row 104 "Torrent"
# and everything below it down by one row), using the same Institution1/
# Address1 (DIPC, San Sebastian/Donostia) pairing already used by several
# other DIPC-affiliated authors (e.g. row 18 "Bayo, A.").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=LastName, B=Initial, C=Footnote, D=ListOrder,
#          E=Institution1, F=Address1, G=Institution2, H=Address2
$newRow = 104

# Shift row 104..end down by inserting a fresh row above the current row 104.
$ws.Rows.Item($newRow).Insert() | Out-Null

# Grab the DIPC institution/address text from the row that was just pushed
# down to 105 (previously row 104, "Torrent"), so the new row matches it
# exactly (same shared-string content/formatting as existing DIPC authors).
$institution1 = $ws.Cells.Item($newRow + 1, 5).Value2
$address1 = $ws.Cells.Item($newRow + 1, 6).Value2

$ws.Cells.Item($newRow, 1).Value = "Tonnel\'e"
$ws.Cells.Item($newRow, 2).Value = "C."
$ws.Cells.Item($newRow, 5).Value = $institution1
$ws.Cells.Item($newRow, 6).Value = $address1

# Restore view/selection to roughly where the author was working when they
# added the new row.
$ws.Activate() | Out-Null
$ws.Range("D104").Select() | Out-Null
